# Update view-count (column F) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 24
$wsExpo.Range("F4").Value  = 238
$wsExpo.Range("F6").Value  = 1155
$wsExpo.Range("F17").Value = 1275
$wsExpo.Range("F19").Value = 269
$wsExpo.Range("F20").Value = 1559
$wsExpo.Range("F21").Value = 1312
$wsExpo.Range("F22").Value = 756
$wsExpo.Range("F27").Value = 376
$wsExpo.Range("F28").Value = 3313

# Sheet "本地生活" (Local Life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 776

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 24
$wsAll.Range("F3").Value  = 776
$wsAll.Range("F7").Value  = 238
$wsAll.Range("F10").Value = 1155
$wsAll.Range("F29").Value = 1275
$wsAll.Range("F31").Value = 270
$wsAll.Range("F32").Value = 1559
$wsAll.Range("F33").Value = 1312
$wsAll.Range("F34").Value = 756
$wsAll.Range("F41").Value = 376
$wsAll.Range("F42").Value = 3313
